$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- Create rows 12-14 (new rows) by copying the formatting of row 11 first,
#     so the new cells get the same border/alignment style (s=1) used
#     throughout the rest of the table instead of a bare default style. ---
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Header row ---
$ws.Range("A1").Value = "序号"
$ws.Range("B1").Value = "PIN"
$ws.Range("C1").Value = "方向"
$ws.Range("D1").Value = "说明"
$ws.Range("E1").Value = "丝印"
$ws.Range("F1").Value = "可使用的用途"
$ws.Range("G1").Value = "临时功能"

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "PA6"
$ws.Range("C2").Value = "AI"
$ws.Range("D2").Value = "雷达中频输入,ADC1_CH6"
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = $null

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "PB6"
$ws.Range("C3").Value = "DO"
$ws.Range("D3").Value = "ENRF"
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = "TIMA_4_PWM1" + $nl + "TIMA_6_PWM8"
$ws.Range("G3").Value = $null

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "PB5"
$ws.Range("C4").Value = "PWM"
$ws.Range("D4").Value = "FSK控制输出,TIMERA3_CH2"
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = "TIMA_3_PWM2" + $nl + "TIMA_6_PWM7"
$ws.Range("G4").Value = $null

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "PA1"
$ws.Range("C5").Value = "DO"
$ws.Range("D5").Value = "板载红色LED，低亮高灭"
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = $null

# --- Row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "PA5"
$ws.Range("C6").Value = "DO"
$ws.Range("D6").Value = "板载绿色LED，低亮高灭"
$ws.Range("E6").Value = $null
$ws.Range("F6").Value = $null
$ws.Range("G6").Value = $null

# --- Row 7 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "PH2"
$ws.Range("C7").Value = "TXD"
$ws.Range("D7").Value = "M4_USART3 tx"
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = $null

# --- Row 8 ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "PC13"
$ws.Range("C8").Value = "RXD"
$ws.Range("D8").Value = "M4_USART3 rx"
$ws.Range("E8").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = $null

# --- Row 9 ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "PA13"
$ws.Range("C9").Value = "DIO"
$ws.Range("D9").Value = "SWDIO"
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = $null
$ws.Range("G9").Value = $null

# --- Row 10 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "PA14"
$ws.Range("C10").Value = "DI"
$ws.Range("D10").Value = "SWCLK"
$ws.Range("E10").Value = $null
$ws.Range("F10").Value = $null
$ws.Range("G10").Value = $null

# --- Row 11 ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "PA8"
$ws.Range("C11").Value = "PWM"
$ws.Range("D11").Value = "按用户要求输出,TIMA_1_PWM1"
$ws.Range("E11").Value = "O"
$ws.Range("F11").Value = "TIM6_1_PWMA" + $nl + "TIMA_1_PWM1"
$ws.Range("G11").Value = $null

# --- Row 12 (new) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "PA7"
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = $null
$ws.Range("E12").Value = "O1"
$ws.Range("F12").Value = "TIM6_1_PWMB " + $nl + "TIMA_1_PWM5 " + $nl + "TIMA_3_PWM2"
$ws.Range("G12").Value = $null

# --- Row 13 (new) ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "PA0"
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = $null
$ws.Range("E13").Value = "ADC1"
$ws.Range("F13").Value = "TIMA_2_PWM1"
$ws.Range("G13").Value = $null

# --- Row 14 (new) ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "PA4"
$ws.Range("C14").Value = $null
$ws.Range("D14").Value = $null
$ws.Range("E14").Value = "ADC2"
$ws.Range("F14").Value = "TIMA_3_PWM5"
$ws.Range("G14").Value = $null

# --- Column F widened to fit the multi-line "可使用的用途" text ---
$ws.Columns("F").ColumnWidth = 13.86

# --- Wrap text + left/center alignment for the multi-line "可使用的用途" cells ---
$ws.Range("F3").WrapText = $true
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").VerticalAlignment = -4108

$ws.Range("F4").WrapText = $true
$ws.Range("F4").HorizontalAlignment = -4131
$ws.Range("F4").VerticalAlignment = -4108

$ws.Range("F11").WrapText = $true
$ws.Range("F11").HorizontalAlignment = -4131
$ws.Range("F11").VerticalAlignment = -4108

$ws.Range("F12").WrapText = $true
$ws.Range("F12").HorizontalAlignment = -4131
$ws.Range("F12").VerticalAlignment = -4108

$ws.Range("F13").WrapText = $true
$ws.Range("F13").HorizontalAlignment = -4131
$ws.Range("F13").VerticalAlignment = -4108

$ws.Range("F14").WrapText = $true
$ws.Range("F14").HorizontalAlignment = -4131
$ws.Range("F14").VerticalAlignment = -4108

# --- Row heights to match the wrapped content (2 lines = 27pt, 3 lines = 40.5pt) ---
$ws.Rows(3).RowHeight = 27
$ws.Rows(4).RowHeight = 27
$ws.Rows(11).RowHeight = 27
$ws.Rows(12).RowHeight = 40.5

# --- Selection / view state ---
$ws.Range("E16").Select()
